# Disable alert dialogs (e.g. delete-sheet confirmation)
$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Source sheets as they exist before the edit
$origData = $wb.Worksheets.Item("Original Data")
$procData = $wb.Worksheets.Item("Processed Data")

# ---------------------------------------------------------------------
# 1) Duplicate "Original Data" -> becomes "IR data input"
# ---------------------------------------------------------------------
$origData.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$irInput = $wb.Worksheets.Item($wb.Worksheets.Count)
$irInput.Name = "IR data input"

# ---------------------------------------------------------------------
# 2) New sheet "Credit data input" right after "IR data input",
#    containing only the header row of "Original Data"
# ---------------------------------------------------------------------
$creditInput = $wb.Worksheets.Add($null, $irInput)
$creditInput.Name = "Credit data input"

for ($col = 2; $col -le 12; $col++) {
    $headerValue = $origData.Cells.Item(1, $col).Value()
    $creditInput.Cells.Item(1, $col).Value = $headerValue
}

$creditInputHeader = $creditInput.Range("B1:L1")
$creditInputHeader.Font.Bold = $true
$creditInputHeader.HorizontalAlignment = -4108
$creditInputHeader.VerticalAlignment = -4160
$creditInputHeader.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3) Duplicate "Processed Data" -> becomes "IR Processed Data"
# ---------------------------------------------------------------------
$procData.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$irProcessed = $wb.Worksheets.Item($wb.Worksheets.Count)
$irProcessed.Name = "IR Processed Data"

# ---------------------------------------------------------------------
# 4) New sheet "Credit Processed Data" right after "IR Processed Data",
#    containing only the header row of "Processed Data"
# ---------------------------------------------------------------------
$creditProcessed = $wb.Worksheets.Add($null, $irProcessed)
$creditProcessed.Name = "Credit Processed Data"

for ($col = 2; $col -le 9; $col++) {
    $headerValue = $procData.Cells.Item(1, $col).Value()
    $creditProcessed.Cells.Item(1, $col).Value = $headerValue
}

$creditProcessedHeader = $creditProcessed.Range("B1:I1")
$creditProcessedHeader.Font.Bold = $true
$creditProcessedHeader.HorizontalAlignment = -4108
$creditProcessedHeader.VerticalAlignment = -4160
$creditProcessedHeader.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 5) Remove the original "Processed Data" sheet - its role is now
#    filled by "IR data input" / "IR Processed Data"
# ---------------------------------------------------------------------
$procData.Delete() | Out-Null

# Restore original active sheet
$origData.Activate()
